$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.212.40"
$ws.Range("E2").Value = "  +1.66%  "
$ws.Range("D3").Value = "2.587.99"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "523.02"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.29"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.564"
$ws.Range("E8").Value = "  -0.84%  "
$ws.Range("D9").Value = "2.597.63"
$ws.Range("E9").Value = "  -0.89%  "
$ws.Range("E10").Value = "  -2.12%  "
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("E12").Value = "  -2.25%  "
$ws.Range("E13").Value = "  +2.65%  "
$ws.Range("D14").Value = "3.046.94"
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("D15").Value = "59.116.05"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("E16").Value = "  -0.27%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000133"
$ws.Range("E17").Value = "  -1.06%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.556.82"
$ws.Range("E18").Value = "  -1.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "341.47"
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("E21").Value = "  -2.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.46"
$ws.Range("E22").Value = "  +1.12%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.64"
$ws.Range("E24").Value = "  +1.99%  "
$ws.Range("E25").Value = "  +0.70%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  +0.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.05"
$ws.Range("E28").Value = "  +0.40%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  -3.59%  "
$ws.Range("E31").Value = "  -5.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.58"
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.70"
$ws.Range("E33").Value = "  -0.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.31"
$ws.Range("E34").Value = "  -0.36%  "
$ws.Range("E35").Value = "  -1.78%  "
$ws.Range("E36").Value = "  -2.13%  "
$ws.Range("E37").Value = "  +2.02%  "
$ws.Range("E38").Value = "  +0.37%  "
$ws.Range("E39").Value = "  -4.41%  "
$ws.Range("E40").Value = "  -6.80%  "
$ws.Range("E41").Value = "  -1.01%  "
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "272.00"
$ws.Range("E43").Value = "  -0.56%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.79"
$ws.Range("E45").Value = "  +0.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0950"
$ws.Range("E46").Value = "  -0.99%  "
$ws.Range("E47").Value = "  -1.77%  "
$ws.Range("E48").Value = "  -2.69%  "
$ws.Range("D49").Value = "1.968.10"
$ws.Range("E49").Value = "  -0.42%  "
$ws.Range("E50").Value = "  -0.24%  "
$ws.Range("E51").Value = "  -2.88%  "
